# SAO_SEBASTIAO_DO_CAI.xlsx update
# - Rename "Paineis DARQ" -> "PAINEIS DARQ"
# - Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
# - Remove the "Desarquivamentos Pendentes" sheet entirely

$wb = $excel.ActiveWorkbook

# Deleting a sheet normally pops a confirmation dialog in real Excel;
# suppress it so the automation can proceed unattended.
$excel.DisplayAlerts = $false

# Remove the obsolete "Desarquivamentos Pendentes" tab and all of its data.
$wsPendentes = $wb.Worksheets.Item("Desarquivamentos Pendentes")
$wsPendentes.Delete()

# Rename the remaining tabs to their updated (uppercase / accented) titles.
$wsPaineis = $wb.Worksheets.Item("Paineis DARQ")
$wsPaineis.Name = "PAINEIS DARQ"

$wsRecolhimento = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$wsRecolhimento.Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Keep the originally-selected tab ("PAINEIS DARQ") active/selected.
$wsPaineis.Activate()
